$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.499.99"
$ws.Range("E2").Value = "  +3.65%  "
$ws.Range("D3").Value = "2.421.39"
$ws.Range("E3").Value = "  +2.51%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'313.71"
$ws.Range("E5").Value = "  +3.84%  "
$ws.Range("D6").Value = "'101.32"
$ws.Range("E6").Value = "  +5.81%  "
$ws.Range("D7").Value = "'0.510"
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("D9").Value = "'0.513"
$ws.Range("E9").Value = "  +5.05%  "
$ws.Range("D10").Value = "'35.23"
$ws.Range("E10").Value = "  +3.33%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.0800"
$ws.Range("E11").Value = "  +2.09%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.126"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").Value = "'18.93"
$ws.Range("E13").Value = "  +3.39%  "
$ws.Range("E14").Value = "  +3.13%  "
$ws.Range("D15").Value = "2.800.44"
$ws.Range("D16").Value = "2.398.16"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("D17").Value = "'0.833"
$ws.Range("E17").Value = "  +5.32%  "
$ws.Range("D18").Value = "44.369.11"
$ws.Range("E18").Value = "  +3.41%  "
$ws.Range("D19").Value = "'12.44"
$ws.Range("E19").Value = "  +5.33%  "
$ws.Range("E20").Value = "  +2.13%  "
$ws.Range("D21").Value = "0.0₃0925"
$ws.Range("E21").Value = "  +4.80%  "
$ws.Range("D22").Value = "'68.71"
$ws.Range("E22").Value = "  +1.27%  "
$ws.Range("D23").Value = "'241.62"
$ws.Range("E23").Value = "  +2.86%  "
$ws.Range("E24").Value = "  +5.90%  "
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D27").Value = "'25.15"
$ws.Range("E27").Value = "  +2.85%  "
$ws.Range("E28").Value = "  -4.51%  "
$ws.Range("E29").Value = "  +3.27%  "
$ws.Range("D30").Value = "'33.39"
$ws.Range("E30").Value = "  +4.87%  "
$ws.Range("D31").Value = "'48.51"
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("D33").Value = "'19.31"
$ws.Range("E33").Value = "  +11.85%  "
$ws.Range("D34").Value = "'5.17"
$ws.Range("E34").Value = "  +3.37%  "
$ws.Range("D35").Value = "'0.0772"
$ws.Range("E35").Value = "  +8.71%  "
$ws.Range("E37").Value = "  +4.51%  "
$ws.Range("D38").Value = "'1.89"
$ws.Range("E38").Value = "  +2.56%  "
$ws.Range("D39").Value = "'2.87"
$ws.Range("E39").Value = "  +2.64%  "
$ws.Range("D40").Value = "'122.57"
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("E42").Value = "  +1.58%  "
$ws.Range("D43").Value = "'20.91"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("E44").Value = "  +4.52%  "
$ws.Range("D45").Value = "1.950.80"
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("E46").Value = "  +2.08%  "
$ws.Range("E47").Value = "  +8.05%  "
$ws.Range("D48").Value = "'9.49"
$ws.Range("E48").Value = "  +3.56%  "
$ws.Range("E49").Value = "  +9.72%  "
$ws.Range("D50").Value = "'55.22"
$ws.Range("E50").Value = "  +7.27%  "
$ws.Range("D51").Value = "'74.14"
$ws.Range("E51").Value = "  +4.39%  "
